# Refactor the "synthetic array" of status indicator values:
#   black square "⬛" + "noir"  -> blue book  "📘" + "bleu"
#   red square   "🟥"           -> red book   "📕"
#   orange square "🟧"          -> orange book "📙"
#   green square "🟩"           -> green book "📗"
#
# These values live in the shared string table and are referenced from
# column A (pictogram) and column B (color label) across the worksheet.
# We use Find/Replace over every worksheet's used range so every
# occurrence (shared string) is updated consistently.

$wb = $excel.ActiveWorkbook

function Replace-AllOccurrences($OldValue, $NewValue) {
    foreach ($ws in $wb.Worksheets) {
        $used = $ws.UsedRange
        $null = $used.Replace($OldValue, $NewValue, 2, 1, $false, $false, $false)
    }
}

Replace-AllOccurrences "⬛" "📘"
Replace-AllOccurrences "🟥" "📕"
Replace-AllOccurrences "🟧" "📙"
Replace-AllOccurrences "🟩" "📗"
Replace-AllOccurrences "noir" "bleu"
